# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the Fenrir_Profits workbook sheets
# as described by the source diff (scheduled market-data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1803.6666
$ws.Range("I111").Value = 2188.3333
$ws.Range("J111").Value = 1611.3334
$ws.Range("K111").Value = 6564.999899999999
$ws.Range("L111").Value = 4834.0002
$ws.Range("M111").Value = -3497.999899999999
$ws.Range("N111").Value = -10968.0002
$ws.Range("H115").Value = 1321.4286
$ws.Range("I115").Value = 625
$ws.Range("J115").Value = 5500
$ws.Range("K115").Value = 1875
$ws.Range("L115").Value = 16500
$ws.Range("M115").Value = -308
$ws.Range("N115").Value = -19634
$ws.Range("H126").Value = 34990
$ws.Range("J126").Value = 34990
$ws.Range("L126").Value = 34990
$ws.Range("N126").Value = -44870
$ws.Range("H132").Value = 20919152
$ws.Range("I132").Value = 21543404
$ws.Range("K132").Value = 64630212
$ws.Range("M132").Value = -64627682
$ws.Range("H139").Value = 26082.223
$ws.Range("J139").Value = 26082.223
$ws.Range("L139").Value = 26082.223
$ws.Range("N139").Value = -36362.223

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3477.425
$ws.Range("I61").Value = 3516.1353
$ws.Range("K61").Value = 3516.1353
$ws.Range("M61").Value = -3304.1353
$ws.Range("H136").Value = 3477.425
$ws.Range("I136").Value = 3516.1353
$ws.Range("K136").Value = 10548.4059
$ws.Range("M136").Value = -7998.4059
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 99980
$ws.Range("J42").Value = 99980
$ws.Range("L42").Value = 99980
$ws.Range("N42").Value = -100636
$ws.Range("H59").Value = 40779.668
$ws.Range("J59").Value = 40779.668
$ws.Range("L59").Value = 40779.668
$ws.Range("N59").Value = -42473.668
$ws.Range("H105").Value = 1720.5714
$ws.Range("I105").Value = 1500
$ws.Range("J105").Value = 1757.3334
$ws.Range("K105").Value = 1500
$ws.Range("L105").Value = 1757.3334
$ws.Range("M105").Value = 247
$ws.Range("N105").Value = -5251.3334
$ws.Range("H137").Value = 46250
$ws.Range("J137").Value = 46250
$ws.Range("L137").Value = 46250
$ws.Range("N137").Value = -56450
$ws.Range("H138").Value = 34000
$ws.Range("J138").Value = 34000
$ws.Range("L138").Value = 34000
$ws.Range("N138").Value = -44280
$ws.Range("H140").Value = 32000
$ws.Range("J140").Value = 32000
$ws.Range("L140").Value = 32000
$ws.Range("N140").Value = -42360

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3791901.8
$ws.Range("I58").Value = 5755087
$ws.Range("J58").Value = 16545.691
$ws.Range("K58").Value = 5755087
$ws.Range("L58").Value = 16545.691
$ws.Range("M58").Value = -5754884
$ws.Range("N58").Value = -16951.691
$ws.Range("H107").Value = 280.8611
$ws.Range("I107").Value = 225
$ws.Range("K107").Value = 225
$ws.Range("M107").Value = 1695
$ws.Range("H122").Value = 5716479
$ws.Range("I122").Value = 20409618
$ws.Range("J122").Value = 2480.2222
$ws.Range("K122").Value = 61228854
$ws.Range("L122").Value = 7440.6666
$ws.Range("M122").Value = -61226404
$ws.Range("N122").Value = -12340.6666
$ws.Range("H136").Value = 3791901.8
$ws.Range("I136").Value = 5755087
$ws.Range("J136").Value = 16545.691
$ws.Range("K136").Value = 17265261
$ws.Range("L136").Value = 49637.073
$ws.Range("M136").Value = -17262711
$ws.Range("N136").Value = -54737.073
$ws.Range("H137").Value = 20000
$ws.Range("J137").Value = 20000
$ws.Range("L137").Value = 20000
$ws.Range("N137").Value = -30200
$ws.Range("H138").Value = 36000
$ws.Range("J138").Value = 36000
$ws.Range("L138").Value = 36000
$ws.Range("N138").Value = -46280
$ws.Range("H139").Value = 38784.285
$ws.Range("J139").Value = 38784.285
$ws.Range("L139").Value = 38784.285
$ws.Range("N139").Value = -49064.285
$ws.Range("H140").Value = 32496.666
$ws.Range("J140").Value = 32496.666
$ws.Range("L140").Value = 32496.666
$ws.Range("N140").Value = -42856.666
$ws.Range("H141").Value = 39946.89
$ws.Range("I141").Value = 10000
$ws.Range("J141").Value = 43690.25
$ws.Range("K141").Value = 10000
$ws.Range("L141").Value = 43690.25
$ws.Range("M141").Value = -4820
$ws.Range("N141").Value = -54050.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 512.4666999999999
$ws.Range("I5").Value = 491.92856
$ws.Range("J5").Value = 800
$ws.Range("K5").Value = 1475.78568
$ws.Range("L5").Value = 2400
$ws.Range("M5").Value = -1363.78568
$ws.Range("N5").Value = -2624
$ws.Range("H80").Value = 3386.875
$ws.Range("I80").Value = 1500.5
$ws.Range("J80").Value = 3764.15
$ws.Range("K80").Value = 4501.5
$ws.Range("L80").Value = 11292.45
$ws.Range("M80").Value = -3565.5
$ws.Range("N80").Value = -13164.45
$ws.Range("H83").Value = 3386.875
$ws.Range("I83").Value = 1500.5
$ws.Range("J83").Value = 3764.15
$ws.Range("K83").Value = 13504.5
$ws.Range("L83").Value = 33877.35
$ws.Range("M83").Value = -8824.5
$ws.Range("N83").Value = -43237.35
$ws.Range("H129").Value = 1634.6086
$ws.Range("J129").Value = 1948.1111
$ws.Range("L129").Value = 5844.3333
$ws.Range("N129").Value = -15844.3333
$ws.Range("H135").Value = 512.4666999999999
$ws.Range("I135").Value = 491.92856
$ws.Range("J135").Value = 800
$ws.Range("K135").Value = 4427.35704
$ws.Range("L135").Value = 7200
$ws.Range("M135").Value = -1892.35704
$ws.Range("N135").Value = -12270

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 30780
$ws.Range("J137").Value = 30780
$ws.Range("L137").Value = 30780
$ws.Range("N137").Value = -40980
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 16000
$ws.Range("J139").Value = 16000
$ws.Range("L139").Value = 16000
$ws.Range("N139").Value = -26280
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 35000
$ws.Range("J141").Value = 35000
$ws.Range("L141").Value = 35000
$ws.Range("N141").Value = -45360

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10529331
$ws.Range("I132").Value = 18184118
$ws.Range("J132").Value = 3997.875
$ws.Range("K132").Value = 54552354
$ws.Range("L132").Value = 11993.625
$ws.Range("M132").Value = -54549824
$ws.Range("N132").Value = -17053.625
$ws.Range("H136").Value = 3988.9348
$ws.Range("I136").Value = 5086.484
$ws.Range("J136").Value = 1720.6666
$ws.Range("K136").Value = 15259.452
$ws.Range("L136").Value = 5161.9998
$ws.Range("M136").Value = -12709.452
$ws.Range("N136").Value = -10261.9998
$ws.Range("H138").Value = 58866.668
$ws.Range("J138").Value = 58866.668
$ws.Range("L138").Value = 58866.668
$ws.Range("N138").Value = -69146.66800000001
$ws.Range("H139").Value = 25123.75
$ws.Range("I139").Value = 24500
$ws.Range("J139").Value = 25331.666
$ws.Range("K139").Value = 24500
$ws.Range("L139").Value = 25331.666
$ws.Range("M139").Value = -19360
$ws.Range("N139").Value = -35611.666
$ws.Range("H140").Value = 77976.336
$ws.Range("J140").Value = 77976.336
$ws.Range("L140").Value = 77976.336
$ws.Range("N140").Value = -88336.336
$ws.Range("H141").Value = 59702.5
$ws.Range("J141").Value = 59702.5
$ws.Range("L141").Value = 59702.5
$ws.Range("N141").Value = -70062.5
